$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 36, pushing the old rows 36-38 down to 38-40.
$ws.Rows(36).Insert()
$ws.Rows(36).Insert()

# --- Row 35: update in place (date + price/volume figures) ---
$ws.Cells.Item(35, 4).Value = 45267
$ws.Cells.Item(35, 10).Value = 100
$ws.Cells.Item(35, 11).Value = 7000
$ws.Cells.Item(35, 12).Value = 8000
$ws.Cells.Item(35, 13).Value = 7500
$ws.Cells.Item(35, 16).Value = 750

# --- Row 36: new row (copy of the surrounding pattern, new values) ---
$ws.Cells.Item(36, 1).Value = 1
$ws.Cells.Item(36, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(36, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(36, 4).Value = 45267
$ws.Cells.Item(36, 5).Value = 15
$ws.Cells.Item(36, 6).Value = 100112043
$ws.Cells.Item(36, 7).Value = "Pepino dulce"
$ws.Cells.Item(36, 8).Value = "Cultivar XV región"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 130
$ws.Cells.Item(36, 11).Value = 5000
$ws.Cells.Item(36, 12).Value = 6000
$ws.Cells.Item(36, 13).Value = 5500
$ws.Cells.Item(36, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(36, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(36, 16).Value = 550
$ws.Cells.Item(36, 17).Value = 10
$ws.Cells.Item(36, 18).Value = "Hortaliza"

# --- Row 37: new row ---
$ws.Cells.Item(37, 1).Value = 1
$ws.Cells.Item(37, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(37, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(37, 4).Value = 44221
$ws.Cells.Item(37, 5).Value = 15
$ws.Cells.Item(37, 6).Value = 100112043
$ws.Cells.Item(37, 7).Value = "Pepino dulce"
$ws.Cells.Item(37, 8).Value = "Cultivar XV región"
$ws.Cells.Item(37, 9).Value = "Primera"
$ws.Cells.Item(37, 10).Value = 140
$ws.Cells.Item(37, 11).Value = 5000
$ws.Cells.Item(37, 12).Value = 6000
$ws.Cells.Item(37, 13).Value = 5500
$ws.Cells.Item(37, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(37, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(37, 16).Value = 550
$ws.Cells.Item(37, 17).Value = 10
$ws.Cells.Item(37, 18).Value = "Hortaliza"

# Ensure the date column keeps its date number format for the new rows
$ws.Range("D36:D37").NumberFormat = $ws.Range("D35").NumberFormat
